$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 106.14286
$ws.Range("I55").Value = 36.25
$ws.Range("J55").Value = 199.33333
$ws.Range("K55").Value = 36.25
$ws.Range("L55").Value = 199.33333
$ws.Range("M55").Value = 177.75
$ws.Range("N55").Value = -627.3333299999999
$ws.Range("H108").Value = 99789.234
$ws.Range("J108").Value = 99789.234
$ws.Range("L108").Value = 99789.234
$ws.Range("N108").Value = -107469.234
$ws.Range("H114").Value = 99810.39999999999
$ws.Range("J114").Value = 99810.39999999999
$ws.Range("L114").Value = 99810.39999999999
$ws.Range("N114").Value = -108488.4
$ws.Range("H117").Value = 90019.75
$ws.Range("J117").Value = 90019.75
$ws.Range("L117").Value = 90019.75
$ws.Range("N117").Value = -99197.75
$ws.Range("H118").Value = 620.6
$ws.Range("I118").Value = 669.1111
$ws.Range("K118").Value = 2007.3333
$ws.Range("M118").Value = -350.3332999999998
$ws.Range("H120").Value = 44594.5
$ws.Range("J120").Value = 44594.5
$ws.Range("L120").Value = 44594.5
$ws.Range("N120").Value = -54270.5
$ws.Range("H133").Value = 71148.47
$ws.Range("J133").Value = 71148.47
$ws.Range("L133").Value = 71148.47
$ws.Range("N133").Value = -81268.47
$ws.Range("H134").Value = 99995
$ws.Range("J134").Value = 99995
$ws.Range("L134").Value = 99995
$ws.Range("N134").Value = -110135
$ws.Range("H136").Value = 76659.2
$ws.Range("J136").Value = 76659.2
$ws.Range("L136").Value = 76659.2
$ws.Range("N136").Value = -86859.2
$ws.Range("H138").Value = 2312.682
$ws.Range("J138").Value = 2479.6365
$ws.Range("L138").Value = 7438.9095
$ws.Range("N138").Value = -17718.9095
$ws.Range("H139").Value = 98402
$ws.Range("J139").Value = 98402
$ws.Range("L139").Value = 98402
$ws.Range("N139").Value = -108682
$ws.Range("H140").Value = 80765.2
$ws.Range("J140").Value = 80765.2
$ws.Range("L140").Value = 80765.2
$ws.Range("N140").Value = -91125.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 89496
$ws.Range("J7").Value = 89496
$ws.Range("L7").Value = 89496
$ws.Range("N7").Value = -89724
$ws.Range("H96").Value = 25000
$ws.Range("J96").Value = 25000
$ws.Range("L96").Value = 25000
$ws.Range("N96").Value = -30492
$ws.Range("H127").Value = 94996.664
$ws.Range("J127").Value = 94996.664
$ws.Range("L127").Value = 94996.664
$ws.Range("N127").Value = -104916.664
$ws.Range("H132").Value = 2147.5625
$ws.Range("I132").Value = 2003.24
$ws.Range("K132").Value = 6009.72
$ws.Range("M132").Value = -3479.72

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1332.1111
$ws.Range("J20").Value = 1407.6666
$ws.Range("L20").Value = 1407.6666
$ws.Range("N20").Value = -1901.6666
$ws.Range("H52").Value = 99990
$ws.Range("J52").Value = 99990
$ws.Range("L52").Value = 99990
$ws.Range("N52").Value = -100516
$ws.Range("H94").Value = 1908.9615
$ws.Range("I94").Value = 1901.591
$ws.Range("K94").Value = 1901.591
$ws.Range("M94").Value = -1450.591
$ws.Range("H95").Value = 40624
$ws.Range("J95").Value = 40624
$ws.Range("L95").Value = 40624
$ws.Range("N95").Value = -46116
$ws.Range("H108").Value = 96281.42999999999
$ws.Range("J108").Value = 96281.42999999999
$ws.Range("L108").Value = 96281.42999999999
$ws.Range("N108").Value = -103961.43
$ws.Range("H110").Value = 81172.5
$ws.Range("J110").Value = 81172.5
$ws.Range("L110").Value = 81172.5
$ws.Range("N110").Value = -89352.5
$ws.Range("H114").Value = 89662.5
$ws.Range("J114").Value = 89662.5
$ws.Range("L114").Value = 89662.5
$ws.Range("N114").Value = -98340.5
$ws.Range("H115").Value = 76996.836
$ws.Range("J115").Value = 79996
$ws.Range("L115").Value = 79996
$ws.Range("N115").Value = -83130
$ws.Range("H116").Value = 66745
$ws.Range("J116").Value = 66745
$ws.Range("L116").Value = 66745
$ws.Range("N116").Value = -75923
$ws.Range("H121").Value = 99990
$ws.Range("J121").Value = 99990
$ws.Range("L121").Value = 99990
$ws.Range("N121").Value = -103484
$ws.Range("H122").Value = 71963.42999999999
$ws.Range("J122").Value = 71963.42999999999
$ws.Range("L122").Value = 71963.42999999999
$ws.Range("N122").Value = -81763.42999999999
$ws.Range("H127").Value = 61267.43
$ws.Range("J127").Value = 61267.43
$ws.Range("L127").Value = 61267.43
$ws.Range("N127").Value = -71187.42999999999
$ws.Range("H132").Value = 29961.309
$ws.Range("J132").Value = 29961.309
$ws.Range("L132").Value = 29961.309
$ws.Range("N132").Value = -40081.309
$ws.Range("H140").Value = 43500
$ws.Range("J140").Value = 43500
$ws.Range("L140").Value = 43500
$ws.Range("N140").Value = -53860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 53891
$ws.Range("J18").Value = 53891
$ws.Range("L18").Value = 53891
$ws.Range("N18").Value = -54351
$ws.Range("H22").Value = 1290.1875
$ws.Range("I22").Value = 1391
$ws.Range("K22").Value = 1391
$ws.Range("M22").Value = -1041
$ws.Range("H116").Value = 82723
$ws.Range("J116").Value = 82723
$ws.Range("L116").Value = 82723
$ws.Range("N116").Value = -91901
$ws.Range("H117").Value = 38898.4
$ws.Range("J117").Value = 38898.4
$ws.Range("L117").Value = 38898.4
$ws.Range("N117").Value = -48076.4
$ws.Range("H132").Value = 1750748.2
$ws.Range("I132").Value = 1978737.6
$ws.Range("J132").Value = 2830
$ws.Range("K132").Value = 5936212.800000001
$ws.Range("L132").Value = 8490
$ws.Range("M132").Value = -5933682.800000001
$ws.Range("N132").Value = -13550
$ws.Range("H138").Value = 54897.5
$ws.Range("J138").Value = 54897.5
$ws.Range("L138").Value = 54897.5
$ws.Range("N138").Value = -65177.5
$ws.Range("H140").Value = 60709
$ws.Range("I140").Value = 60709
$ws.Range("K140").Value = 60709
$ws.Range("M140").Value = -55529

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 75000024
$ws.Range("I29").Value = 99
$ws.Range("K29").Value = 297
$ws.Range("M29").Value = -20
$ws.Range("H107").Value = 872.1667
$ws.Range("J107").Value = 777
$ws.Range("L107").Value = 2331
$ws.Range("N107").Value = -6171

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 99995
$ws.Range("J110").Value = 99995
$ws.Range("L110").Value = 99995
$ws.Range("N110").Value = -108175
$ws.Range("H113").Value = 3300
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830
$ws.Range("H116").Value = 59996.57
$ws.Range("J116").Value = 59996.57
$ws.Range("L116").Value = 59996.57
$ws.Range("N116").Value = -69174.57000000001
$ws.Range("H119").Value = 51763.383
$ws.Range("J119").Value = 51763.383
$ws.Range("L119").Value = 51763.383
$ws.Range("N119").Value = -61439.383
$ws.Range("H132").Value = 3866.0605
$ws.Range("I132").Value = 3110
$ws.Range("K132").Value = 9330
$ws.Range("M132").Value = -6800
$ws.Range("H135").Value = 44785.5
$ws.Range("J135").Value = 44785.5
$ws.Range("L135").Value = 44785.5
$ws.Range("N135").Value = -54925.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2803.7222
$ws.Range("I16").Value = 2386.1333
$ws.Range("J16").Value = 4891.6665
$ws.Range("K16").Value = 2386.1333
$ws.Range("L16").Value = 4891.6665
$ws.Range("M16").Value = -2216.1333
$ws.Range("N16").Value = -5231.6665
$ws.Range("H22").Value = 817
$ws.Range("I22").Value = 805.44446
$ws.Range("K22").Value = 805.44446
$ws.Range("M22").Value = -510.44446
$ws.Range("H27").Value = 817
$ws.Range("I27").Value = 805.44446
$ws.Range("K27").Value = 805.44446
$ws.Range("M27").Value = -698.44446
$ws.Range("H117").Value = 82580
$ws.Range("J117").Value = 89096
$ws.Range("L117").Value = 89096
$ws.Range("N117").Value = -98274
$ws.Range("H123").Value = 75667.27
$ws.Range("J123").Value = 78195
$ws.Range("L123").Value = 78195
$ws.Range("N123").Value = -87995
$ws.Range("H132").Value = 2319.1333
$ws.Range("H136").Value = 3823.2104
$ws.Range("J136").Value = 2426.818
$ws.Range("L136").Value = 7280.454000000001
$ws.Range("N136").Value = -12380.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 36997.332
$ws.Range("J121").Value = 36997.332
$ws.Range("L121").Value = 36997.332
$ws.Range("N121").Value = -40491.332
$ws.Range("H122").Value = 2102.9285
$ws.Range("I122").Value = 1290.6364
$ws.Range("K122").Value = 3871.9092
$ws.Range("M122").Value = -1421.9092
